$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 18020
$ws.Range("B3").Value = 38136
$ws.Range("B4").Value = 55101
$ws.Range("B5").Value = 24223
$ws.Range("B6").Value = 24637
$ws.Range("B7").Value = 25596
$ws.Range("B8").Value = 87754
$ws.Range("B9").Value = 15843
$ws.Range("B10").Value = 22657
$ws.Range("B11").Value = 28983
$ws.Range("B12").Value = 185187
$ws.Range("B13").Value = 47158
$ws.Range("B14").Value = 127244
$ws.Range("B15").Value = 32827
$ws.Range("B16").Value = 42337
$ws.Range("B17").Value = 154815
$ws.Range("B18").Value = 23155
$ws.Range("B19").Value = 41652
$ws.Range("B20").Value = 11307
$ws.Range("B21").Value = 16984
$ws.Range("B22").Value = 39970
$ws.Range("B23").Value = 34668
$ws.Range("B24").Value = 18883
$ws.Range("B25").Value = 34571
$ws.Range("B26").Value = 1792
$ws.Range("B27").Value = 35879
$ws.Range("B28").Value = 32124
$ws.Range("B29").Value = 11692
$ws.Range("B30").Value = 13787
$ws.Range("B31").Value = 36614
$ws.Range("B32").Value = 17784
$ws.Range("B33").Value = 89780
$ws.Range("B34").Value = 40622
$ws.Range("B35").Value = 38486
$ws.Range("B36").Value = 15654
$ws.Range("B37").Value = 50896
$ws.Range("B38").Value = 42943
$ws.Range("B39").Value = 14779
$ws.Range("B40").Value = 46538
$ws.Range("B41").Value = 37012
$ws.Range("B42").Value = 23334
$ws.Range("B43").Value = 28105
$ws.Range("B44").Value = 8874
$ws.Range("B45").Value = 35180
$ws.Range("B46").Value = 23693
$ws.Range("B47").Value = 17781
$ws.Range("B48").Value = 39492
$ws.Range("B49").Value = 4755
$ws.Range("B50").Value = 10589
$ws.Range("B51").Value = 73077
$ws.Range("B52").Value = 37460
$ws.Range("B53").Value = 28727
$ws.Range("B54").Value = 23386
$ws.Range("B55").Value = 29241
$ws.Range("B56").Value = 21434
$ws.Range("B57").Value = 15985
$ws.Range("B58").Value = 20399
$ws.Range("B59").Value = 137491
$ws.Range("B60").Value = 23638
$ws.Range("B61").Value = 15187
$ws.Range("B62").Value = 47060
$ws.Range("B63").Value = 22729
$ws.Range("B64").Value = 25789
$ws.Range("B65").Value = 31879
$ws.Range("B66").Value = 342259
$ws.Range("B67").Value = 37693
$ws.Range("B68").Value = 39345
$ws.Range("B69").Value = 20281
$ws.Range("B70").Value = 66852
$ws.Range("B71").Value = 75293
$ws.Range("B72").Value = 67740
$ws.Range("B73").Value = 28287
$ws.Range("B74").Value = 35395
$ws.Range("B75").Value = 54414
$ws.Range("B76").Value = 8249
$ws.Range("B77").Value = 22093
$ws.Range("B78").Value = 30426
$ws.Range("B79").Value = 34172
